$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 10016.2
$ws.Range("I8").Value = 20.25
$ws.Range("J8").Value = 50000
$ws.Range("K8").Value = 60.75
$ws.Range("L8").Value = 150000
$ws.Range("M8").Value = 78.25
$ws.Range("N8").Value = -150278

$ws.Range("H18").Value = 1785.7142
$ws.Range("I18").Value = 1783.3334
$ws.Range("J18").Value = 1800
$ws.Range("K18").Value = 1783.3334
$ws.Range("L18").Value = 1800
$ws.Range("M18").Value = -1499.3334
$ws.Range("N18").Value = -2368

$ws.Range("H19").Value = 720
$ws.Range("I19").Value = 700
$ws.Range("J19").Value = 800
$ws.Range("K19").Value = 700
$ws.Range("L19").Value = 800
$ws.Range("M19").Value = -525
$ws.Range("N19").Value = -1150

$ws.Range("H55").Value = 365
$ws.Range("I55").Value = 407.66666
$ws.Range("J55").Value = 301
$ws.Range("K55").Value = 407.66666
$ws.Range("L55").Value = 301
$ws.Range("M55").Value = -193.66666
$ws.Range("N55").Value = -729

$ws.Range("H76").Value = 4320
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()

$ws.Range("H79").Value = 4320
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

$ws.Range("H107").Value = 540
$ws.Range("I107").Value = 572.94116
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 572.94116
$ws.Range("L107").Value = 400
$ws.Range("M107").Value = 1347.05884
$ws.Range("N107").Value = -4240

$ws.Range("H109").Value = 60249.715
$ws.Range("J109").Value = 60249.715
$ws.Range("L109").Value = 60249.715
$ws.Range("N109").Value = -63023.715

$ws.Range("H129").Value = 1004.9375
$ws.Range("I129").Value = 622.5
$ws.Range("J129").Value = 1059.5714
$ws.Range("K129").Value = 1867.5
$ws.Range("L129").Value = 3178.7142
$ws.Range("M129").Value = 3132.5
$ws.Range("N129").Value = -13178.7142

$ws.Range("H138").Value = 1680.24
$ws.Range("I138").Value = 1023.3774
$ws.Range("J138").Value = 2420.9575
$ws.Range("K138").Value = 3070.1322
$ws.Range("L138").Value = 7262.872499999999
$ws.Range("M138").Value = 2069.8678
$ws.Range("N138").Value = -17542.8725

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 500006750
$ws.Range("I43").Value = 15000
$ws.Range("J43").Value = 666670660
$ws.Range("K43").Value = 15000
$ws.Range("L43").Value = 666670660
$ws.Range("M43").Value = -14687
$ws.Range("N43").Value = -666671286

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1447.7222
$ws.Range("I94").Value = 1178.6666
$ws.Range("J94").Value = 1985.8334
$ws.Range("K94").Value = 1178.6666
$ws.Range("L94").Value = 1985.8334
$ws.Range("M94").Value = -727.6666
$ws.Range("N94").Value = -2887.8334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 598.34784
$ws.Range("I107").Value = 496.66666
$ws.Range("J107").Value = 613.6
$ws.Range("K107").Value = 496.66666
$ws.Range("L107").Value = 613.6
$ws.Range("M107").Value = 1423.33334
$ws.Range("N107").Value = -4453.6

$ws.Range("H134").Value = 1696.7241
$ws.Range("I134").Value = 1482.619
$ws.Range("J134").Value = 2258.75
$ws.Range("K134").Value = 4447.857
$ws.Range("L134").Value = 6776.25
$ws.Range("M134").Value = -1912.857
$ws.Range("N134").Value = -11846.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 1000001
$ws.Range("J9").Value = 1000001
$ws.Range("L9").Value = 3000003
$ws.Range("N9").Value = -3000451

$ws.Range("H136").Value = 3610.2727
$ws.Range("I136").Value = 1060
$ws.Range("J136").Value = 4177
$ws.Range("K136").Value = 3180
$ws.Range("L136").Value = 12531
$ws.Range("M136").Value = 1920
$ws.Range("N136").Value = -22731

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 20250
$ws.Range("J93").Value = 20250
$ws.Range("L93").Value = 20250
$ws.Range("N93").Value = -23994

$ws.Range("H122").Value = 4973.1816
$ws.Range("J122").Value = 3400
$ws.Range("L122").Value = 10200
$ws.Range("N122").Value = -15100

$ws.Range("H123").Value = 18889.533
$ws.Range("J123").Value = 18889.533
$ws.Range("L123").Value = 18889.533
$ws.Range("N123").Value = -23789.533

$ws.Range("H132").Value = 1954.9756
$ws.Range("I132").Value = 1482.7646
$ws.Range("J132").Value = 4248.5713
$ws.Range("K132").Value = 4448.293799999999
$ws.Range("L132").Value = 12745.7139
$ws.Range("M132").Value = -1918.293799999999
$ws.Range("N132").Value = -17805.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 51500
$ws.Range("J15").Value = 100000
$ws.Range("L15").Value = 100000
$ws.Range("N15").Value = -100576

$ws.Range("H20").Value = 30005.5
$ws.Range("J20").Value = 10011
$ws.Range("L20").Value = 10011
$ws.Range("N20").Value = -10491

$ws.Range("H21").Value = 50000
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws.Range("H25").Value = 24256.75
$ws.Range("J25").Value = 24256.75
$ws.Range("L25").Value = 24256.75
$ws.Range("N25").Value = -24842.75

$ws.Range("H35").Value = 50000
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

$ws.Range("H41").Value = 500003000
$ws.Range("J41").Value = 500003000
$ws.Range("L41").Value = 500003000
$ws.Range("N41").Value = -500003780

$ws.Range("H45").Value = 22333.334
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 22333.334
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 22333.334
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -23315.334

$ws.Range("H74").Value = 6280.6665
$ws.Range("I74").Value = 3400
$ws.Range("J74").Value = 7103.7144
$ws.Range("K74").Value = 3400
$ws.Range("L74").Value = 7103.7144
$ws.Range("M74").Value = -2464
$ws.Range("N74").Value = -8975.714400000001

$ws.Range("H77").Value = 6280.6665
$ws.Range("I77").Value = 3400
$ws.Range("J77").Value = 7103.7144
$ws.Range("K77").Value = 10200
$ws.Range("L77").Value = 21311.1432
$ws.Range("M77").Value = -5520
$ws.Range("N77").Value = -30671.1432

$ws.Range("H92").Value = 30509.8
$ws.Range("J92").Value = 30509.8
$ws.Range("L92").Value = 30509.8
$ws.Range("N92").Value = -35501.8

$ws.Range("H122").Value = 11365681
$ws.Range("I122").Value = 13159147
$ws.Range("J122").Value = 7065
$ws.Range("K122").Value = 39477441
$ws.Range("L122").Value = 21195
$ws.Range("M122").Value = -39474991
$ws.Range("N122").Value = -26095

$ws.Range("H124").Value = 60607
$ws.Range("J124").Value = 60607
$ws.Range("L124").Value = 60607
$ws.Range("N124").Value = -70427
